$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report (2023-03-23, serial 45008) needs to be inserted
# right after the existing row for serial 44257 (row 398), pushing the
# following historical rows down by two and adding the new "Primera" /
# "Segunda" quality rows at the top of that block.
$ws.Rows("399:400").Insert()

$ws.Range("A399").Value = 9
$ws.Range("B399").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C399").Value = "Metropolitana"
$ws.Range("D399").Value = 45008
$ws.Range("E399").Value = 13
$ws.Range("F399").Value = 100112017
$ws.Range("G399").Value = "Apio"
$ws.Range("H399").Value = "Americana (o)"
$ws.Range("I399").Value = "Primera"
$ws.Range("J399").Value = 70
$ws.Range("K399").Value = 8000
$ws.Range("L399").Value = 9000
$ws.Range("M399").Value = 8500
$ws.Range("N399").Value = "$/docena de matas"
$ws.Range("O399").Value = "Región de Coquimbo"
$ws.Range("P399").Value = 1417
$ws.Range("Q399").Value = 6
$ws.Range("R399").Value = "Hortaliza"

$ws.Range("A400").Value = 9
$ws.Range("B400").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C400").Value = "Metropolitana"
$ws.Range("D400").Value = 45008
$ws.Range("E400").Value = 13
$ws.Range("F400").Value = 100112017
$ws.Range("G400").Value = "Apio"
$ws.Range("H400").Value = "Americana (o)"
$ws.Range("I400").Value = "Segunda"
$ws.Range("J400").Value = 43
$ws.Range("K400").Value = 7000
$ws.Range("L400").Value = 7000
$ws.Range("M400").Value = 7000
$ws.Range("N400").Value = "$/docena de matas"
$ws.Range("O400").Value = "Región de Coquimbo"
$ws.Range("P400").Value = 1167
$ws.Range("Q400").Value = 6
$ws.Range("R400").Value = "Hortaliza"
